$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '56.059.70'
$ws.Range("E2").Value = '  -1.39%  '
$ws.Range("D3").Value = '2.977.34'
$ws.Range("E3").Value = '  -0.09%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '503.11'
$ws.Range("E5").Value = '  +0.95%  '
$ws.Range("D6").Value = '137.16'
$ws.Range("E6").Value = '  -0.31%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.426'
$ws.Range("E8").Value = '  -0.74%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.15'
$ws.Range("E9").Value = '  -1.38%  '
$ws.Range("E10").Value = '  -0.82%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.364'
$ws.Range("E11").Value = '  +1.85%  '
$ws.Range("D12").Value = '3.481.15'
$ws.Range("E12").Value = '  -0.42%  '
$ws.Range("E13").Value = '  -1.34%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.87'
$ws.Range("E14").Value = '  -0.58%  '
$ws.Range("E15").Value = '  +0.84%  '
$ws.Range("D16").Value = '56.040.25'
$ws.Range("E16").Value = '  -1.64%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.00'
$ws.Range("E17").Value = '  -0.85%  '
$ws.Range("D18").Value = '2.973.10'
$ws.Range("E18").Value = '  -0.53%  '
$ws.Range("D19").Value = '12.82'
$ws.Range("E19").Value = '  +1.98%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.95'
$ws.Range("E20").Value = '  +1.34%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '326.79'
$ws.Range("E21").Value = '  +2.29%  '
$ws.Range("D22").Value = '0.999'
$ws.Range("E22").Value = '  -0.03%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.493'
$ws.Range("E23").Value = '  +0.93%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '64.48'
$ws.Range("E24").Value = '  +1.50%  '
$ws.Range("D25").Value = '3.094.44'
$ws.Range("E25").Value = '  -0.41%  '
$ws.Range("E26").Value = '  -0.19%  '
$ws.Range("E27").Value = '  +0.01%  '
$ws.Range("D28").Value = '0.0₃0915'
$ws.Range("E28").Value = '  +3.04%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.37'
$ws.Range("E29").Value = '  -2.16%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.98'
$ws.Range("E30").Value = '  -1.03%  '
$ws.Range("E31").Value = '  +0.49%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.16'
$ws.Range("E32").Value = '  -0.06%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '20.14'
$ws.Range("E33").Value = '  -0.21%  '
$ws.Range("D34").Value = '152.87'
$ws.Range("E34").Value = '  -1.54%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.49'
$ws.Range("E35").Value = '  -1.41%  '
$ws.Range("D36").Value = '5.75'
$ws.Range("E36").Value = '  -0.18%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '25.80'
$ws.Range("E37").Value = '  +6.37%  '
$ws.Range("E38").Value = '  -0.46%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0659'
$ws.Range("E39").Value = '  -0.65%  '
$ws.Range("D40").Value = '3.007.02'
$ws.Range("E40").Value = '  -0.17%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '36.89'
$ws.Range("E41").Value = '  -2.10%  '
$ws.Range("E42").Value = '  -0.05%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.78'
$ws.Range("E43").Value = '  +1.22%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.649'
$ws.Range("E44").Value = '  +1.20%  '
$ws.Range("D45").Value = '2.168.22'
$ws.Range("E45").Value = '  -1.49%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.35'
$ws.Range("E46").Value = '  -2.42%  '
$ws.Range("B47").Value = 'ONDO'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.920'
$ws.Range("E47").Value = '  -1.71%  '
$ws.Range("B48").Value = 'Cosmos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.81'
$ws.Range("E48").Value = '  -2.13%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0236'
$ws.Range("E49").Value = '  +0.82%  '
$ws.Range("D50").Value = '19.48'
$ws.Range("E50").Value = '  +1.36%  '
$ws.Range("B51").Value = 'Notcoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/2L2Y4ghjj+notcoin-not'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0167'
$ws.Range("E51").Value = '  +13.12%  '
